$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the shift labels in column A (A, B, C -> Morning, Afternoon, Night)
$ws.Range("A1").Value = "Morning"
$ws.Range("A2").Value = "Afternoon"
$ws.Range("A3").Value = "Night"

# Move the active selection to A3
$ws.Range("A3").Select()
